$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in FACTORYSERIAL value within the "Pasos a seguir" text of row 4 (F4)
$row4Height = $ws.Rows.Item(4).RowHeight
$f4 = $ws.Range("F4").Text
$f4New = $f4.Replace("485724435AXXXXXX", "4.5724435AXXXXXX")
$ws.Range("F4").Value = $f4New
$ws.Rows.Item(4).RowHeight = $row4Height

# Add new "Pasos a seguir" text for row 5 (F5), which was previously empty
$ws.Range("F5").Value = "1.Seleccionar la barra de b" + [char]0x00FA + "squeda y digitar " + [char]0x0022 + "HUAWEI TEST" + [char]0x0022 + [char]10 + "2."

# Update the active selection to F5 to match the saved workbook state
$ws.Range("F5").Select()
